# "remove shape spline for double-blind review"
#
# The resume's publication list contains a citation for:
#   "Q. Xia, C. Chen, S. Li*, A. Hao and H. Qin. Fast 4D Shape Sequence
#    Completion from Sparse Samples via Spline Fitting in Linear Rotation
#    Invariant Space. Graphical Models. (CCF B, under review)"
# This paragraph must be removed in its entirety (together with the
# paragraph mark) so double-blind reviewers cannot identify the author
# through a self-citation to an in-review paper.

$d = $word.ActiveDocument

# Locate the paragraph holding the citation by its distinctive text
# rather than a hard-coded index, searching from the end of the
# document since the citation is the final bullet in the publications
# list.
$target = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*Fast 4D Shape Sequence Completion*Spline Fitting*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    # Deleting the paragraph's Range also removes its trailing paragraph
    # mark, merging what remains (nothing) into the preceding paragraph -
    # exactly mirroring the two <w:p> elements collapsing into one in the
    # underlying OOXML.
    $target.Range.Delete()
}
